$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "plain" decimal number (e.g. "0.470", "479.40")
# must be pinned to Text format first, otherwise Excel auto-converts the
# typed string to a number and silently drops the significant trailing zero
# (e.g. "0.470" -> 0.47), which would not match the source data.
$textValueCells = @(
    "D5",
    "D6",
    "D11",
    "D12",
    "D13",
    "D20",
    "D21",
    "D24",
    "D25",
    "D28",
    "D29",
    "D30",
    "D33",
    "D34",
    "D35",
    "D38",
    "D40",
    "D41",
    "D47",
    "D48",
)
foreach ($addr in $textValueCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.544.98"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "3.160.08"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "613.92"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "145.12"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.156.52"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("D11").Value = "5.41"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").Value = "0.470"
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").Value = "3.681.45"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("E16").Value = "  +3.60%  "
$ws.Range("D17").Value = "64.572.52"
$ws.Range("E17").Value = "  +0.76%  "
$ws.Range("D18").Value = "3.159.37"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").Value = "479.40"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("D21").Value = "14.58"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("E22").Value = "  +2.59%  "
$ws.Range("E23").Value = "  +3.49%  "
$ws.Range("D24").Value = "13.77"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").Value = "83.80"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -3.42%  "
$ws.Range("D28").Value = "8.63"
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("D29").Value = "7.12"
$ws.Range("E29").Value = "  +5.72%  "
$ws.Range("D30").Value = "0.119"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("E31").Value = "  -5.20%  "
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "26.47"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").Value = "2.67"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("D35").Value = "1.12"
$ws.Range("E35").Value = "  +2.15%  "
$ws.Range("D36").Value = "0.0₃0784"
$ws.Range("E36").Value = "  +8.03%  "
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").Value = "53.19"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("E39").Value = "  +3.91%  "
$ws.Range("D40").Value = "460.98"
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("D41").Value = "0.0398"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").Value = "2.859.34"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("E45").Value = "  +3.06%  "
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "2.46"
$ws.Range("E47").Value = "  +6.72%  "
$ws.Range("D48").Value = "26.51"
$ws.Range("E48").Value = "  +0.56%  "
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("E50").Value = "  +9.25%  "
$ws.Range("E51").Value = "  -0.41%  "
